# "Mise a jour de l'application" - add the next training/match session
# (16 Sep 2025) as a new column (AU) on the attendance sheet, with each
# player's status ("P" = present, "B" = blessure) for that date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> attendance letter for the new session (column AU).
$attendance = @{
    2  = "P"; 3  = "P"; 4  = "P"; 5  = "B"; 6  = "B";
    7  = "P"; 8  = "P"; 9  = "P"; 10 = "B"; 11 = "P";
    12 = "P"; 13 = "B"; 14 = "P"; 15 = "P"; 16 = "P";
    17 = "P"; 18 = "P"; 19 = "P"; 20 = "P"; 21 = "B";
    22 = "P"; 23 = "P"; 24 = "P"; 25 = "P"; 26 = "P";
    27 = "P"; 28 = "P"; 29 = "P"
}

# New date header for column AU (16 Sep 2025 -> Excel serial 45916).
# Write the value first (so dependent formulas pick up the change), then
# copy the neighbouring header cell's format only, so the existing date
# style (s="5") is reused instead of a brand new style being created.
$ws.Range("AU1").Value = 45916
$ws.Range("AT1").Copy()
$ws.Range("AU1").PasteSpecial(-4122)

foreach ($row in 2..29) {
    $letter = $attendance[$row]
    $cell = $ws.Range("AU$row")
    $cell.Value = $letter
    $ws.Range("AT$row").Copy()
    $cell.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# Match the author's final selection.
$ws.Range("AW25").Select()
